# Refresh the cryptos price/volume snapshot (coinranking.com scrape),
# commit "Updated cryptos list on Mon May  8 08:57:49 UTC 2023 with GitHub Actions".
#
# Column D ("Price") cells are forced to Text number-format ("@") right before
# the write: several of the new price strings (e.g. "1.008", "317.26") would
# otherwise be auto-parsed by Excel into numeric values, silently dropping
# trailing zeros / precision and losing the original "NNN.NNN" textual layout
# used by this sheet (prices with thousand separators like "27.922.12").
# Column E ("Volume(1h)") values already contain '%' and padding spaces, so
# Excel keeps them as plain text without any extra formatting step.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.922.12"
$ws.Range("E2").Value = "  -3.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.854.57"
$ws.Range("E3").Value = "  -2.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.26"
$ws.Range("E5").Value = "  -2.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4360"
$ws.Range("E7").Value = "  -5.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3672"
$ws.Range("E8").Value = "  -3.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07483"
$ws.Range("E9").Value = "  -3.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9357"
$ws.Range("E10").Value = "  -4.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.29"
$ws.Range("E11").Value = "  -3.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.858.01"
$ws.Range("E12").Value = "  -2.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.683"
$ws.Range("E13").Value = "  -3.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.404"
$ws.Range("E14").Value = "  -4.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06889"
$ws.Range("E15").Value = "  -2.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.007"
$ws.Range("E16").Value = "  +0.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.20"
$ws.Range("E17").Value = "  -3.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008980"
$ws.Range("E18").Value = "  -5.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.85"
$ws.Range("E20").Value = "  -5.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.932.82"
$ws.Range("E21").Value = "  -3.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.101"
$ws.Range("E22").Value = "  -4.25%  "

$ws.Range("E23").Value = "  -0.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.095.92"
$ws.Range("E24").Value = "  -2.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.013"
$ws.Range("E25").Value = "  -3.72%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.26"
$ws.Range("E26").Value = "  -2.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.33"
$ws.Range("E27").Value = "  -3.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.352"
$ws.Range("E28").Value = "  -5.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.96"
$ws.Range("E29").Value = "  -3.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.728"
$ws.Range("E30").Value = "  -6.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08969"
$ws.Range("E31").Value = "  -3.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7974"
$ws.Range("E32").Value = "  -8.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.816"
$ws.Range("E33").Value = "  -5.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.003"
$ws.Range("E34").Value = "  -2.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.168"
$ws.Range("E35").Value = "  -6.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.005"
$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.119"
$ws.Range("E37").Value = "  -3.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05407"
$ws.Range("E38").Value = "  -5.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01951"
$ws.Range("E39").Value = "  -4.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.921"
$ws.Range("E40").Value = "  +2.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5228"
$ws.Range("E41").Value = "  -4.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.978"
$ws.Range("E42").Value = "  -5.63%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1676"
$ws.Range("E43").Value = "  -4.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.708"
$ws.Range("E44").Value = "  -6.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06699"
$ws.Range("E45").Value = "  -2.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4862"
$ws.Range("E46").Value = "  -6.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.53"
$ws.Range("E47").Value = "  -5.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.42"
$ws.Range("E48").Value = "  -3.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.005"
$ws.Range("E49").Value = "  +0.25%  "

$ws.Range("E50").Value = "  -7.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.668"
$ws.Range("E51").Value = "  -6.30%  "
